$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1) - column F ("想去人数") updates
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 24
$wsExpo.Range("F4").Value = 13
$wsExpo.Range("F5").Value = 3919
$wsExpo.Range("F6").Value = 164
$wsExpo.Range("F8").Value = 239
$wsExpo.Range("F9").Value = 19

# Sheet "全部类型" (index 4) - column F ("想去人数") updates
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 24
$wsAll.Range("F8").Value = 13
$wsAll.Range("F9").Value = 3921
$wsAll.Range("F10").Value = 164
$wsAll.Range("F13").Value = 239
$wsAll.Range("F14").Value = 19
